# Quarterly indexing esoteric bug-fix: shift each revision date from the
# 1st of its quarter-start month to the 15th of the following month.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 150; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $serial = $cell.Value2
    $d = Get-Date -Year 1899 -Month 12 -Day 30 -Hour 0 -Minute 0 -Second 0
    $d = $d.AddDays($serial)
    $next = $d.AddMonths(1)
    $fixed = Get-Date -Year $next.Year -Month $next.Month -Day 15 -Hour 0 -Minute 0 -Second 0
    $cell.Value2 = $fixed
}
